# The deck's color theme ("Integral") is replaced by the default
# "Office Theme" color palette. All slides share the single slide
# master's theme, so re-pointing the theme colours through any slide's
# ThemeColorScheme updates the shared theme part used by the whole
# presentation.
#
# Theme colour slot order (ThemeColorSchemeIndex 1-12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
#
# Target "Office Theme" palette (RRGGBB), converted below to the
# PowerPoint RGB() long (0x00BBGGRR) expected by ThemeColor.RGB:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

$p = $ppt.ActivePresentation

$officeThemeColors = @{
    1  = 0x000000   # dk1      000000
    2  = 0xFFFFFF   # lt1      FFFFFF
    3  = 0x6A5444   # dk2      44546A
    4  = 0xE6E6E7   # lt2      E7E6E6
    5  = 0xD59B5B   # accent1  5B9BD5
    6  = 0x317DED   # accent2  ED7D31
    7  = 0xA5A5A5   # accent3  A5A5A5
    8  = 0x00C0FF   # accent4  FFC000
    9  = 0xC47244   # accent5  4472C4
    10 = 0x47AD70   # accent6  70AD47
    11 = 0xC16305   # hlink    0563C1
    12 = 0x724F95   # folHlink 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i]
}
